# Append a new "2021年" data row (row 11) to Sheet1, mirroring the
# existing yearly rows (2012年..2020年) already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

# Column A: the year label, styled like the other year cells (e.g. A10).
$ws.Range("A$row").Value = "2021年"
$ws.Range("A10").Copy()
$ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats

# Column E has no data for this year (matches the blank pattern already
# present for E6:E10 - those cells hold an empty string, which Excel
# treats the same as a cleared/blank cell, so it is simply left unset).

# Remaining columns: numeric economic indicators for 2021.
$values = @{
    "B"  = 1895.76
    "C"  = 581.85
    "D"  = 96.87
    "F"  = 909.5599999999999
    "G"  = 5566.03
    "H"  = 612.52
    "I"  = 5608.91
    "J"  = 230.66
    "K"  = 66715.00999999999
    "L"  = 250.89
    "M"  = 93.09
    "N"  = -52.49
    "O"  = 446.21
    "P"  = 2940.17
    "Q"  = 406.23
    "R"  = 126.58
    "S"  = 924.85
    "T"  = 173.59
    "U"  = 4901.28
    "V"  = 1167.93
    "W"  = 5932.79
    "X"  = 312.3
    "Y"  = 174.5
    "Z"  = 3618.35
    "AA" = 377.64
    "AB" = 2584.92
    "AC" = 1465.42
    "AD" = 769.47
    "AE" = 534.78
    "AF" = 7773.04
    "AG" = 2026.21
    "AH" = 712.9
    "AI" = 2258.8
    "AJ" = 38.87
    "AK" = 1384.77
    "AL" = 518.64
    "AM" = 3210.15
    "AN" = 80.63
    "AO" = 1211.82
    "AP" = 4250.3
    "AQ" = 597.96
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$row").Value = $values[$col]
}
